$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The original sheet had a two-row header (row1 + row2). Collapse it into a
# single header row by deleting the old second header row; this shifts the
# three data rows up from 3/4/5 to 2/3/4.
$ws.Rows("2").Delete()

# Rewrite row 1 as the new single header row.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:E1 carry no special formatting (default style) - some of these cells
# held leftover formatting from the old header row, so reset the font back
# to the sheet default to drop it.
$ws.Range("A1:E1").Font.Size = 10

# F1:K1 use the small (9pt) header font, like the rest of the data cells.
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").VerticalAlignment = -4107

# Update the active selection to point at the first data row, matching the
# post-edit workbook state.
$ws.Range("A2:K2").Select()
